# Apply the data updates described by the commit diff ("remove_25" — one
# study/observation removed from the underlying meta-analysis, shifting the
# mean/min/max/o/n summary figures on all three sheets).

$wb = $excel.ActiveWorkbook

$wsMean = $wb.Worksheets.Item("mean_effect")
$wsTWW  = $wb.Worksheets.Item("TWW")
$wsSoil = $wb.Worksheets.Item("Soil")

# ---------------------------------------------------------------------
# mean_effect sheet
# ---------------------------------------------------------------------

# row 4 - pH
$wsMean.Range("E4").Value = 1
$wsMean.Range("F4").Value = 39
$wsMean.Range("G4").Value = 12
$wsMean.Range("I4").Value = -2

# row 7 - porosity
$wsMean.Range("C7").Value = -6
$wsMean.Range("D7").Value = -10
$wsMean.Range("E7").Value = -2
$wsMean.Range("F7").Value = 18
$wsMean.Range("G7").Value = 6
$wsMean.Range("I7").Value = -4

# row 8 - bulk density
$wsMean.Range("C8").Value = -1
$wsMean.Range("D8").Value = -3
$wsMean.Range("F8").Value = 42
$wsMean.Range("G8").Value = 12
$wsMean.Range("I8").Value = -1

# ---------------------------------------------------------------------
# TWW sheet
# ---------------------------------------------------------------------

# row 15 - pH Mean
$wsTWW.Range("F15").Value = 1
$wsTWW.Range("G15").Value = 39
$wsTWW.Range("H15").Value = 12

# row 30 - p Mean
$wsTWW.Range("D30").Value = -6
$wsTWW.Range("E30").Value = -10
$wsTWW.Range("F30").Value = -2
$wsTWW.Range("G30").Value = 18
$wsTWW.Range("H30").Value = 6

# row 35 - bulk density Mean
$wsTWW.Range("D35").Value = -1
$wsTWW.Range("E35").Value = -3
$wsTWW.Range("G35").Value = 42
$wsTWW.Range("H35").Value = 12

# ---------------------------------------------------------------------
# Soil sheet
# ---------------------------------------------------------------------

# row 13 - pH Medium
$wsSoil.Range("F13").Value = 1
$wsSoil.Range("G13").Value = 26
$wsSoil.Range("H13").Value = 8
$wsSoil.Range("J13").Value = -2

# row 14 - pH Coarse
$wsSoil.Range("D14").Value = -22
$wsSoil.Range("E14").Value = -24
$wsSoil.Range("F14").Value = -20
$wsSoil.Range("G14").Value = 2
$wsSoil.Range("H14").Value = 1
$wsSoil.Range("I14").Formula = "=2"
$wsSoil.Range("J14").Value = -2

# row 15 - pH Mean
$wsSoil.Range("F15").Value = 1
$wsSoil.Range("G15").Value = 39
$wsSoil.Range("H15").Value = 12
$wsSoil.Range("J15").Value = -2

# row 28 - p Medium
$wsSoil.Range("D28").Value = -11
$wsSoil.Range("E28").Value = -17
$wsSoil.Range("F28").Value = -4
$wsSoil.Range("G28").Value = 10
$wsSoil.Range("H28").Value = 3
$wsSoil.Range("I28").Value = 6
$wsSoil.Range("J28").Value = -7

# row 29 - p Coarse
$wsSoil.Range("E29").Value = -1
$wsSoil.Range("F29").Value = 6
$wsSoil.Range("G29").Value = 1
$wsSoil.Range("H29").Value = 1
$wsSoil.Range("I29").Value = 3
$wsSoil.Range("J29").Value = -4

# row 30 - p Mean
$wsSoil.Range("D30").Value = -6
$wsSoil.Range("E30").Value = -10
$wsSoil.Range("F30").Value = -2
$wsSoil.Range("G30").Value = 18
$wsSoil.Range("H30").Value = 6
$wsSoil.Range("J30").Value = -4

# row 33 - bulk density Medium
$wsSoil.Range("D33").Value = -1
$wsSoil.Range("E33").Value = -4
$wsSoil.Range("F33").Value = 1
$wsSoil.Range("G33").Value = 27
$wsSoil.Range("H33").Value = 8

# row 34 - bulk density Coarse
$wsSoil.Range("D34").Value = 1
$wsSoil.Range("G34").Value = 1
$wsSoil.Range("H34").Value = 1
$wsSoil.Range("I34").Value = 3
$wsSoil.Range("J34").Value = -2

# row 35 - bulk density Mean
$wsSoil.Range("D35").Value = -1
$wsSoil.Range("E35").Value = -3
$wsSoil.Range("G35").Value = 42
$wsSoil.Range("H35").Value = 12
$wsSoil.Range("J35").Value = -1

# ---------------------------------------------------------------------
# View / selection state (cosmetic - matches the saved workbook's cursor
# position on each sheet and which sheet/tab is active on re-open).
# ---------------------------------------------------------------------

# mean_effect: selection moved from C3:G3 to C7:I7
$wsMean.Range("C7:I7").Select()

# TWW: no longer the active tab, selection collapses to J31
$wsTWW.Range("J31").Select()

# Soil: becomes the active tab, selection collapses to J30
$wsSoil.Activate()
$wsSoil.Range("J30").Select()
